$ws = $excel.ActiveWorkbook.Worksheets.Item("Test Cases")
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 0
$win.SplitColumn = 0
$win.ScrollRow = 17
$win.ScrollColumn = 2
